# Updates cryptos list values (price + volume%) per upstream diff (Oct 25 2023 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.033.98"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "1.790.12"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.12"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.97"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.281"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0661"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "2.047.86"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.36"
$ws.Range("E14").Value = "  +10.41%  "
$ws.Range("D15").Value = "1.790.65"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.635"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "34.063.60"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.60"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "253.03"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").Value = "0.0₃0744"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.43"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  -2.83%  "
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.44"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.60"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.02"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "1.468.46"
$ws.Range("E36").Value = "  -7.70%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.632"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0186"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.62"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "1.946.59"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.72"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.29"
$ws.Range("E51").Value = "  -4.75%  "
